# Fruta / hortaliza, semanal
# Insert a new weekly record at the top of the data table (row 33), shifting
# the existing rows 33-41 down to 34-42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 33, shifting cells down.
$ws.Rows.Item(33).Insert(-4121)   # -4121 = xlShiftDown

# Populate the newly inserted row 33 with the new weekly record.
$ws.Cells.Item(33, 1).Value  = 1
$ws.Cells.Item(33, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(33, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(33, 4).Value  = 44523
$ws.Cells.Item(33, 5).Value  = 15
$ws.Cells.Item(33, 6).Value  = 100112031
$ws.Cells.Item(33, 7).Value  = "Poroto verde"
$ws.Cells.Item(33, 8).Value  = "Sin especificar"
$ws.Cells.Item(33, 9).Value  = "Primera"
$ws.Cells.Item(33, 10).Value = 1300
$ws.Cells.Item(33, 11).Value = 800
$ws.Cells.Item(33, 12).Value = 850
$ws.Cells.Item(33, 13).Value = 825
$ws.Cells.Item(33, 14).Value = "$/kilo"
$ws.Cells.Item(33, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(33, 16).Value = 825
$ws.Cells.Item(33, 17).Value = 1
$ws.Cells.Item(33, 18).Value = "Hortaliza"
